# Polish Jinja2 Template Render
# Rename Sheet1 -> Train, add a new empty "Test" sheet, and make it the active sheet.

$wb = $excel.ActiveWorkbook

$trainSheet = $wb.Worksheets.Item(1)
$trainSheet.Name = "Train"

$testSheet = $wb.Worksheets.Add($null, $trainSheet)
$testSheet.Name = "Test"

$testSheet.Activate()
